# Pytania do Promotora.xlsx - "poprawiono update, czytelnosc zapytan sql"
#
# Net change:
#  - The old "Potestować wszystkie endpointy na wszelkie sposoby" row is
#    removed; the row that used to hold "Potworzyć gotowe zbiory call'i w
#    PostManie" moves up one line (A5) and the row it vacated (A6) becomes
#    a blank separator row.
#  - A new "TO_ASK" block is appended below the blank separator, mirroring
#    the existing "TO_DO" block: a bold/shaded header cell followed by
#    bordered detail rows (one of which wraps onto several lines) and a
#    couple of blank bordered rows at the end.
#  - Selection moves to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the "Potestować..." row: shift "Potworzyć..." up into A5,
#        and blank out what is now the trailing A6 (keeping it as an empty
#        row rather than leaving stale text behind). ---
$ws.Range("A5").Value = "Potworzyć gotowe zbiory call’i w PostManie"
$ws.Range("A6").Clear()

# --- 2. Build the new "TO_ASK" section in rows 7-11, copying formatting
#        from the matching rows of the existing "TO_DO" section so the
#        header keeps the bold/shaded style (s=2) and the detail rows keep
#        the bordered/wrapped style (s=1). ---

# Row 7: section header, same look as A1 ("TO_DO").
$ws.Range("A1").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "TO_ASK"
$ws.Range("A7").RowHeight = 15.75

# Row 8: short question line, same look as A2.
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "GminaRelacja - zmienić na EFF_DT?"
$ws.Range("A8").RowHeight = 15.75

# Row 9: longer explanatory line - taller row to fit the wrapped text.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Z powodu GminaRelacja update kolumny start_date który jest jednocześnie kluczem główny zaczyna wywalać z powodu tabeli GminaRelacja"
$ws.Range("A9").RowHeight = 45.75

# Rows 10-11: trailing blank bordered rows, same look as A2 but no text.
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").RowHeight = 15.75

$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").RowHeight = 15.75

$excel.CutCopyMode = $false

# --- 3. Move the selection to A8, matching the saved view state. ---
$ws.Range("A8").Select()
